$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new data rows beneath the existing "Browser"/"URL" table
$ws.Range("A3").Value = "TS1"
$ws.Range("B3").Value = "Danish"
$ws.Range("A4").Value = "TS2"
$ws.Range("B4").Value = "Check"

# Match the author's final selection (B4)
$ws.Range("B4").Select()
